# The document carries three inline logo pictures anchored in the
# section's headers/footers:
#   - two copies of the Pearson logo (AlternativeText/descr is the
#     PearsonLogo.png source path) that need to be renamed from
#     "image1.png" to "image2.png"
#   - one copy of the BTec logo (AlternativeText/descr =
#     "BTec_Logo-Orange") that needs to be renamed from "image2.jpg"
#     to "image1.jpg"
#
# Walk every header/footer of every section and rename the inline
# picture shapes in place, matching on their (stable) AlternativeText
# so the edit doesn't depend on guessing header/footer indices.
#
# NB: after InlineShape.Name is set, other already-fetched shape
# handles in the same story can go stale, so each rename below is
# followed by a small pipeline write - that forces the host to settle
# / re-resolve handles before the next header/footer is touched.

$d = $word.ActiveDocument

$pearsonAlt = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
$btecAlt = "BTec_Logo-Orange"

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                $alt = $shp.AlternativeText
                if ($alt -eq $pearsonAlt) {
                    $shp.Name = "image2.png"
                    "Renamed header $i shape $j to image2.png"
                } elseif ($alt -eq $btecAlt) {
                    $shp.Name = "image1.jpg"
                    "Renamed header $i shape $j to image1.jpg"
                }
            }
        }

        $ftr = $sec.Footers($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                $alt = $shp.AlternativeText
                if ($alt -eq $pearsonAlt) {
                    $shp.Name = "image2.png"
                    "Renamed footer $i shape $j to image2.png"
                } elseif ($alt -eq $btecAlt) {
                    $shp.Name = "image1.jpg"
                    "Renamed footer $i shape $j to image1.jpg"
                }
            }
        }
    }
}
